# Updates the cryptos price/volume table with refreshed values.
# Column D ("Price") cells are forced to Text via NumberFormat "@" before
# assignment so numeric-looking strings (e.g. "312.47", "1.00", "13.00")
# are not auto-coerced to floating point numbers by Excel, which would
# corrupt precision/trailing zeros. Column E ("Volume(1h)") values already
# contain padding spaces and a "%" sign so they remain text naturally.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.921.55'
$ws.Range('E2').Value = '  +2.35%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.338.27'
$ws.Range('E3').Value = '  +2.05%  '
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.47'
$ws.Range('E5').Value = '  -0.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '108.17'
$ws.Range('E6').Value = '  +2.89%  '
$ws.Range('E7').Value = '  +0.90%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.618'
$ws.Range('E9').Value = '  +2.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.15'
$ws.Range('E10').Value = '  +3.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0917'
$ws.Range('E11').Value = '  +1.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.54'
$ws.Range('E12').Value = '  +1.25%  '
$ws.Range('E13').Value = '  -1.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.00'
$ws.Range('E14').Value = '  +0.47%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.46'
$ws.Range('E15').Value = '  +1.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.695.92'
$ws.Range('E16').Value = '  +2.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.331.11'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.885.02'
$ws.Range('E18').Value = '  +2.46%  '
$ws.Range('E19').Value = '  +1.72%  '
$ws.Range('E20').Value = '  +1.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.00'
$ws.Range('E21').Value = '  -6.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '74.13'
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('E23').Value = '  -2.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '268.55'
$ws.Range('E24').Value = '  +0.96%  '
$ws.Range('E25').Value = '  +2.29%  '
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.63'
$ws.Range('E27').Value = '  +4.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.10'
$ws.Range('E28').Value = '  +2.14%  '
$ws.Range('E29').Value = '  -1.94%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.99'
$ws.Range('E30').Value = '  +4.74%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.59'
$ws.Range('E31').Value = '  +0.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '169.01'
$ws.Range('E32').Value = '  +1.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0886'
$ws.Range('E33').Value = '  +1.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.83'
$ws.Range('E34').Value = '  +9.58%  '
$ws.Range('E35').Value = '  +0.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.114'
$ws.Range('E36').Value = '  +1.28%  '
$ws.Range('E37').Value = '  +3.73%  '
$ws.Range('E38').Value = '  +3.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.91'
$ws.Range('E39').Value = '  +9.14%  '
$ws.Range('E40').Value = '  -0.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.71'
$ws.Range('E41').Value = '  +8.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '104.64'
$ws.Range('E42').Value = '  +10.93%  '
$ws.Range('E43').Value = '  +2.13%  '
$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '71.53'
$ws.Range('E44').Value = '  +0.94%  '
$ws.Range('B45').Value = 'Celestia'
$ws.Range('C45').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.39'
$ws.Range('E45').Value = '  +9.88%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '113.78'
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.676.20'
$ws.Range('E48').Value = '  -3.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '76.81'
$ws.Range('E49').Value = '  -4.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.215'
$ws.Range('E51').Value = '  +13.48%  '
